# Update: po 25. 01. 2021
# Applies the daily COVID stats update for Slovakia to rows 292-326 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously reported AgTests/AgPosit (columns H/I) values ---
$ws.Cells.Item(292, 8).Value = 81215

$ws.Cells.Item(306, 8).Value = 70727

$ws.Cells.Item(309, 8).Value = 57318
$ws.Cells.Item(309, 9).Value = 3964

$ws.Cells.Item(315, 8).Value = 66051
$ws.Cells.Item(315, 9).Value = 3010

$ws.Cells.Item(317, 8).Value = 61022
$ws.Cells.Item(317, 9).Value = 2124

$ws.Cells.Item(318, 8).Value = 24692
$ws.Cells.Item(318, 9).Value = 888

$ws.Cells.Item(320, 8).Value = 86035
$ws.Cells.Item(320, 9).Value = 4194

$ws.Cells.Item(321, 8).Value = 86604
$ws.Cells.Item(321, 9).Value = 2686

$ws.Cells.Item(322, 8).Value = 102224
$ws.Cells.Item(322, 9).Value = 2255

$ws.Cells.Item(323, 8).Value = 143426
$ws.Cells.Item(323, 9).Value = 2235

# --- Row 324 (2021-01-22) corrected values ---
$ws.Cells.Item(324, 2).Value = 234571
$ws.Cells.Item(324, 4).Value = 43923
$ws.Cells.Item(324, 6).Value = 1544
$ws.Cells.Item(324, 8).Value = 209823
$ws.Cells.Item(324, 9).Value = 2458

# --- New row 325 (2021-01-23) ---
$ws.Cells.Item(325, 1).Value = 44219
$ws.Cells.Item(325, 2).Value = 236476
$ws.Cells.Item(325, 3).Value = 189959
$ws.Cells.Item(325, 4).Value = 42449
$ws.Cells.Item(325, 5).Value = 16498
$ws.Cells.Item(325, 6).Value = 1905
$ws.Cells.Item(325, 7).Value = 4068
$ws.Cells.Item(325, 8).Value = 562488
$ws.Cells.Item(325, 9).Value = 5346

# --- New row 326 (2021-01-24) ---
$ws.Cells.Item(326, 1).Value = 44220
$ws.Cells.Item(326, 2).Value = 237027
$ws.Cells.Item(326, 3).Value = 193380
$ws.Cells.Item(326, 4).Value = 39579
$ws.Cells.Item(326, 5).Value = 3806
$ws.Cells.Item(326, 6).Value = 551
$ws.Cells.Item(326, 7).Value = 4068
$ws.Cells.Item(326, 8).Value = 284768
$ws.Cells.Item(326, 9).Value = 2322

# Ensure the date cells in the new rows use the same date number format
# ("yyyy-mm-dd", the same one applied throughout column A) as the rest of
# the column (style is normally inherited automatically from the column
# default, but we set it explicitly to be safe).
$ws.Range("A325").NumberFormat = "yyyy-mm-dd"
$ws.Range("A326").NumberFormat = "yyyy-mm-dd"
